$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve exact formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "42.791.42"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.533.85"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "304.08"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "97.65"
$ws.Range("E6").Value = "  +5.84%  "
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "0.0820"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "7.70"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "2.930.24"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "2.483.74"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "15.02"
$ws.Range("E16").Value = "  +5.73%  "
$ws.Range("D17").Value = "0.868"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "42.823.00"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "71.51"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "253.23"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "2.92"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "2.07"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").Value = "27.69"
$ws.Range("E26").Value = "  -4.59%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  +7.01%  "
$ws.Range("D30").Value = "37.95"
$ws.Range("D31").Value = "6.22"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").Value = "157.15"
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").Value = "19.40"
$ws.Range("E33").Value = "  +13.07%  "
$ws.Range("D34").Value = "2.12"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").Value = "3.30"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "0.0793"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "24.81"
$ws.Range("E39").Value = "  +7.89%  "
$ws.Range("D40").Value = "0.119"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +30.09%  "
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "2.092.03"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "86.26"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "2.786.44"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "73.38"
$ws.Range("E50").Value = "  +6.47%  "
$ws.Range("D51").Value = "0.191"
$ws.Range("E51").Value = "  +1.57%  "
